$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-CellText "D2" "252.01"
Set-CellText "D3" "21.71"
Set-CellText "D4" "5.555"
Set-CellText "D5" "0.05681"
Set-CellText "D6" "6.456"
Set-CellText "D8" "1.043"
Set-CellText "B9" "One"
Set-CellText "C9" "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-CellText "D9" "0.01168"
Set-CellText "E9" "8OneONEBestin24h"
Set-CellText "B10" "WazirX"
Set-CellText "C10" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-CellText "D10" "0.1431"
Set-CellText "E10" "9WazirXWRX"
Set-CellText "B11" "MandalaExchangeToken"
Set-CellText "C11" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-CellText "D11" "0.07313"
Set-CellText "E11" "10MandalaExchangeTokenMDX"
Set-CellText "B12" "LiechtensteinCryptoassetsExchange"
Set-CellText "C12" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-CellText "D12" "0.03148"
Set-CellText "E12" "11LiechtensteinCryptoassetsExchangeLCX"
Set-CellText "B13" "BitrueCoin"
Set-CellText "C13" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-CellText "D13" "0.02940"
Set-CellText "E13" "12BitrueCoinBTR"
Set-CellText "B14" "BitMartToken"
Set-CellText "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-CellText "D14" "0.09263"
Set-CellText "E14" "13BitMartTokenBMX"
Set-CellText "B15" "BitForexToken"
Set-CellText "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-CellText "D15" "0.001664"
Set-CellText "E15" "14BitForexTokenBF"
Set-CellText "B16" "MCDex"
Set-CellText "C16" "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-CellText "D16" "3.224"
Set-CellText "E16" "15MCDexMCB"
Set-CellText "B17" "CoinExToken"
Set-CellText "C17" "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-CellText "D17" "0.04780"
Set-CellText "E17" "16CoinExTokenCET"
Set-CellText "D18" "0.006461"
Set-CellText "D19" "0.005046"
Set-CellText "E19" "18HotbitTokenHTB"
Set-CellText "D20" "0.001054"
Set-CellText "D22" "0.0003203"
Set-CellText "D23" "3.986"
Set-CellText "D24" "3.386"
Set-CellText "D25" "2.089"
Set-CellText "D40" "0.04129"
Set-CellText "D41" "0.006895"
Set-CellText "B42" "CEJI"
Set-CellText "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-CellText "D42" "0.003503"
Set-CellText "E42" "41CEJICEJI"
Set-CellText "B43" "BKEXToken"
Set-CellText "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-CellText "D43" "0.1044"
Set-CellText "E43" "42BKEXTokenBKK"
Set-CellText "D44" "0.009571"
Set-CellText "D45" "0.00005647"
Set-CellText "D46" "0.00000000751"
Set-CellText "D47" "0.7858"
Set-CellText "D48" "0.01703"
Set-CellText "D49" "0.00002102"
Set-CellText "D50" "0.01011"
